$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Max number of SSTs per skid
$ws.Range("B2").Value = 5

# Update Apparent Power at 40 deg C (base value used by formulas below)
$ws.Range("B3").Value = 1000

# B4 and B5 lose their formulas and become plain static values
$ws.Range("B4").Value = 1000
$ws.Range("B5").Value = 1000

# B6 and B7 keep their formulas; they will recalculate automatically
# based on the new B3 value, but ensure formulas remain as specified.
$ws.Range("B6").Formula = "=`$B`$3*(1-0.01)^15"
$ws.Range("B7").Formula = "=`$B`$3*(1-0.01)^20"

# Update the selected cell to B6
$ws.Range("B6").Select()

$wb.Save()
